$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1624.3132
$ws.Range("I9").Value = 435.68832
$ws.Range("J9").Value = 16878.334
$ws.Range("K9").Value = 435.68832
$ws.Range("L9").Value = 16878.334
$ws.Range("M9").Value = -266.68832
$ws.Range("N9").Value = -17216.334

$ws.Range("H28").Value = 913.0526
$ws.Range("I28").Value = 857.38464
$ws.Range("J28").Value = 1033.6666
$ws.Range("K28").Value = 857.38464
$ws.Range("L28").Value = 1033.6666
$ws.Range("M28").Value = -372.38464
$ws.Range("N28").Value = -2003.6666

$ws.Range("H54").Value = 20000
$ws.Range("J54").Value = 20000
$ws.Range("L54").Value = 20000
$ws.Range("N54").Value = -20972

$ws.Range("H64").Value = 4750
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 4750
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H82").Value = 920.25
$ws.Range("I82").Value = 766.3333
$ws.Range("J82").Value = 1382
$ws.Range("K82").Value = 2298.9999
$ws.Range("L82").Value = 4146
$ws.Range("M82").Value = -1892.9999
$ws.Range("N82").Value = -4958

$ws.Range("H85").Value = 920.25
$ws.Range("I85").Value = 766.3333
$ws.Range("J85").Value = 1382
$ws.Range("K85").Value = 2298.9999
$ws.Range("L85").Value = 4146
$ws.Range("M85").Value = -894.9998999999998
$ws.Range("N85").Value = -6954

$ws.Range("H116").Value = 6800.1514
$ws.Range("I116").Value = 6422.2964
$ws.Range("K116").Value = 6422.2964
$ws.Range("M116").Value = -2980.2964

$ws.Range("H133").Value = 126110.57
$ws.Range("J133").Value = 126110.57
$ws.Range("L133").Value = 126110.57
$ws.Range("N133").Value = -136230.57

$ws.Range("H138").Value = 2540.8267
$ws.Range("I138").Value = 1318.8064
$ws.Range("J138").Value = 3401.7954
$ws.Range("K138").Value = 3956.4192
$ws.Range("L138").Value = 10205.3862
$ws.Range("M138").Value = 1183.5808
$ws.Range("N138").Value = -20485.3862

$ws.Range("H139").Value = 71034.75
$ws.Range("J139").Value = 71034.75
$ws.Range("L139").Value = 71034.75
$ws.Range("N139").Value = -81314.75

$ws.Range("H140").Value = 89999.625
$ws.Range("J140").Value = 89999.625
$ws.Range("L140").Value = 89999.625
$ws.Range("N140").Value = -100359.625

$ws.Range("H141").Value = 4948.6665
$ws.Range("I141").Value = 4738.4
$ws.Range("K141").Value = 14215.2
$ws.Range("M141").Value = -9035.199999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 17811.5
$ws.Range("I19").Value = 648.5
$ws.Range("J19").Value = 34974.5
$ws.Range("K19").Value = 648.5
$ws.Range("L19").Value = 34974.5
$ws.Range("M19").Value = -419.5
$ws.Range("N19").Value = -35432.5

$ws.Range("H45").Value = 3407.8333
$ws.Range("I45").Value = 1189.5555
$ws.Range("K45").Value = 1189.5555
$ws.Range("M45").Value = -812.5554999999999

$ws.Range("H61").Value = 38229.184
$ws.Range("J61").Value = 55963.668
$ws.Range("L61").Value = 55963.668
$ws.Range("N61").Value = -56387.668

$ws.Range("H122").Value = 50615.895
$ws.Range("I122").Value = 2670.3845
$ws.Range("K122").Value = 8011.1535
$ws.Range("M122").Value = -5561.1535

$ws.Range("H132").Value = 2881.8276
$ws.Range("I132").Value = 2758.7917
$ws.Range("J132").Value = 3472.4
$ws.Range("K132").Value = 8276.375100000001
$ws.Range("L132").Value = 10417.2
$ws.Range("M132").Value = -5746.375100000001
$ws.Range("N132").Value = -15477.2

$ws.Range("H136").Value = 38229.184
$ws.Range("J136").Value = 55963.668
$ws.Range("L136").Value = 167891.004
$ws.Range("N136").Value = -172991.004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 41799.6
$ws.Range("J58").Value = 41799.6
$ws.Range("L58").Value = 41799.6
$ws.Range("N58").Value = -42387.6

$ws.Range("H105").Value = 3999
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws.Range("H138").Value = 96594.61
$ws.Range("J138").Value = 96594.61
$ws.Range("L138").Value = 96594.61
$ws.Range("N138").Value = -106874.61

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 19232.072
$ws.Range("I4").Value = 19999.75
$ws.Range("J4").Value = 18925
$ws.Range("K4").Value = 19999.75
$ws.Range("L4").Value = 18925
$ws.Range("M4").Value = -19887.75
$ws.Range("N4").Value = -19149

$ws.Range("H31").Value = 1590.9767
$ws.Range("I31").Value = 1187.7646
$ws.Range("J31").Value = 3114.2222
$ws.Range("K31").Value = 1187.7646
$ws.Range("L31").Value = 3114.2222
$ws.Range("M31").Value = -892.7646
$ws.Range("N31").Value = -3704.2222

$ws.Range("H34").Value = 1590.9767
$ws.Range("I34").Value = 1187.7646
$ws.Range("J34").Value = 3114.2222
$ws.Range("K34").Value = 1187.7646
$ws.Range("L34").Value = 3114.2222
$ws.Range("M34").Value = -985.7646
$ws.Range("N34").Value = -3518.2222

$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52080

$ws.Range("H135").Value = 67399.71000000001
$ws.Range("J135").Value = 67399.71000000001
$ws.Range("L135").Value = 67399.71000000001
$ws.Range("N135").Value = -77539.71000000001

$ws.Range("H138").Value = 83653.69
$ws.Range("J138").Value = 83653.69
$ws.Range("L138").Value = 83653.69
$ws.Range("N138").Value = -93933.69

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 23170.855
$ws.Range("I4").Value = 12498.568
$ws.Range("J4").Value = 119221.445
$ws.Range("K4").Value = 37495.704
$ws.Range("L4").Value = 357664.335
$ws.Range("M4").Value = -37383.704
$ws.Range("N4").Value = -357888.335

$ws.Range("H58").Value = 3110.75
$ws.Range("I58").Value = 3147.6667
$ws.Range("K58").Value = 9443.000100000001
$ws.Range("M58").Value = -9315.000100000001

$ws.Range("H62").Value = 8254.666999999999
$ws.Range("J62").Value = 9998.799999999999
$ws.Range("L62").Value = 29996.4
$ws.Range("N62").Value = -31368.4

$ws.Range("H65").Value = 8254.666999999999
$ws.Range("J65").Value = 9998.799999999999
$ws.Range("L65").Value = 89989.2
$ws.Range("N65").Value = -96853.2

$ws.Range("H82").Value = 5332.6665
$ws.Range("J82").Value = 10000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30812

$ws.Range("H85").Value = 5332.6665
$ws.Range("J85").Value = 10000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32808

$ws.Range("H131").Value = 4111.8
$ws.Range("I131").Value = 1624.75
$ws.Range("K131").Value = 4874.25
$ws.Range("M131").Value = 165.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 192786.14
$ws.Range("I5").Value = 250200.8
$ws.Range("J5").Value = 49249.5
$ws.Range("K5").Value = 250200.8
$ws.Range("L5").Value = 49249.5
$ws.Range("M5").Value = -250088.8
$ws.Range("N5").Value = -49473.5

$ws.Range("H18").Value = 4996
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H70").Value = 6849.1665
$ws.Range("I70").Value = 5868
$ws.Range("J70").Value = 7830.3335
$ws.Range("K70").Value = 5868
$ws.Range("L70").Value = 7830.3335
$ws.Range("M70").Value = -5598
$ws.Range("N70").Value = -8370.333500000001

$ws.Range("H73").Value = 6849.1665
$ws.Range("I73").Value = 5868
$ws.Range("J73").Value = 7830.3335
$ws.Range("K73").Value = 5868
$ws.Range("L73").Value = 7830.3335
$ws.Range("M73").Value = -4932
$ws.Range("N73").Value = -9702.333500000001

$ws.Range("H113").Value = 4427.619
$ws.Range("I113").Value = 4640.5557
$ws.Range("K113").Value = 4640.5557
$ws.Range("M113").Value = -2470.5557

$ws.Range("H135").Value = 99306.42999999999
$ws.Range("J135").Value = 99306.42999999999
$ws.Range("L135").Value = 99306.42999999999
$ws.Range("N135").Value = -109446.43

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 91848.5
$ws.Range("J2").Value = 12685.25
$ws.Range("L2").Value = 12685.25
$ws.Range("N2").Value = -12909.25

$ws.Range("H22").Value = 5644.119
$ws.Range("I22").Value = 1722.7
$ws.Range("J22").Value = 6869.5625
$ws.Range("K22").Value = 1722.7
$ws.Range("L22").Value = 6869.5625
$ws.Range("M22").Value = -1427.7
$ws.Range("N22").Value = -7459.5625

$ws.Range("H27").Value = 5644.119
$ws.Range("I27").Value = 1722.7
$ws.Range("J27").Value = 6869.5625
$ws.Range("K27").Value = 1722.7
$ws.Range("L27").Value = 6869.5625
$ws.Range("M27").Value = -1615.7
$ws.Range("N27").Value = -7083.5625

$ws.Range("H133").Value = 89997
$ws.Range("J133").Value = 89997
$ws.Range("L133").Value = 89997
$ws.Range("N133").Value = -95057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 4399.6

$ws.Range("H16").Value = 89999
$ws.Range("J16").Value = 89999
$ws.Range("L16").Value = 89999
$ws.Range("N16").Value = -90583

$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524

$ws.Range("H107").Value = 15000
$ws.Range("I107").Value = 10000
$ws.Range("J107").Value = 20000
$ws.Range("K107").Value = 30000
$ws.Range("L107").Value = 60000
$ws.Range("M107").Value = -28080
$ws.Range("N107").Value = -63840

$ws.Range("H136").Value = 3859.0967
$ws.Range("I136").Value = 3457.111
$ws.Range("K136").Value = 10371.333
$ws.Range("M136").Value = -7821.332999999999

$ws.Range("H139").Value = 102310.445
$ws.Range("J139").Value = 102310.445
$ws.Range("L139").Value = 102310.445
$ws.Range("N139").Value = -112590.445
